# Update dSF (column F) values to reflect repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -4
    6  = -4
    7  = -7
    8  = -12
    11 = -6
    12 = -11
    13 = 1
    17 = 1
    18 = -5
    19 = -1
    23 = 2
    28 = -5
    29 = 4
    34 = 8
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
